# Scheduled-runner refresh of market-price-derived columns (H:N) on the
# Atomos_Profits leve-profitability sheets. Source values (currentAveragePrice*,
# LevePrice*, LeveProfit*) are recomputed from a fresh market-board pull;
# leve metadata columns (A:G) are untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1278.3667
$ws.Range("I98").Value = 1073.2693
$ws.Range("J98").Value = 2611.5
$ws.Range("K98").Value = 1073.2693
$ws.Range("L98").Value = 2611.5
$ws.Range("M98").Value = 424.7307000000001
$ws.Range("N98").Value = -5607.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1278.3667
$ws.Range("I122").Value = 1073.2693
$ws.Range("J122").Value = 2611.5
$ws.Range("K122").Value = 3219.8079
$ws.Range("L122").Value = 7834.5
$ws.Range("M122").Value = -769.8078999999998
$ws.Range("N122").Value = -12734.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4880847
$ws.Range("I132").Value = 5407890
$ws.Range("K132").Value = 16223670
$ws.Range("M132").Value = -16221140

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1197.3334
$ws.Range("I135").Value = 1206.875
$ws.Range("J135").Value = 1178.25
$ws.Range("K135").Value = 10861.875
$ws.Range("L135").Value = 10604.25
$ws.Range("M135").Value = -8326.875
$ws.Range("N135").Value = -15674.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 6257219
$ws.Range("I137").Value = 9099182
$ws.Range("J137").Value = 4900
$ws.Range("K137").Value = 27297546
$ws.Range("L137").Value = 14700
$ws.Range("M137").Value = -27294996
$ws.Range("N137").Value = -19800

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3246.7856
$ws.Range("I138").Value = 1598.4783
$ws.Range("J138").Value = 5242.1055
$ws.Range("K138").Value = 4795.4349
$ws.Range("L138").Value = 15726.3165
$ws.Range("M138").Value = 344.5650999999998
$ws.Range("N138").Value = -26006.3165

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 465201.53
$ws.Range("I141").Value = 1702.8096
$ws.Range("J141").Value = 977489.5600000001
$ws.Range("K141").Value = 5108.4288
$ws.Range("L141").Value = 2932468.68
$ws.Range("M141").Value = 71.57120000000032
$ws.Range("N141").Value = -2942828.68

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 314.75
$ws.Range("I4").Value = 86.333336
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 86.333336
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = 29.666664
$ws.Range("N4").Value = -1232

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 253
$ws.Range("I5").Value = 221.42857
$ws.Range("J5").Value = 326.66666
$ws.Range("K5").Value = 221.42857
$ws.Range("L5").Value = 326.66666
$ws.Range("M5").Value = -109.42857
$ws.Range("N5").Value = -550.66666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1493.8125
$ws.Range("I74").Value = 1562.4546
$ws.Range("J74").Value = 1342.8
$ws.Range("K74").Value = 1562.4546
$ws.Range("L74").Value = 1342.8
$ws.Range("M74").Value = -688.4546
$ws.Range("N74").Value = -3090.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1493.8125
$ws.Range("I77").Value = 1562.4546
$ws.Range("J77").Value = 1342.8
$ws.Range("K77").Value = 7812.273
$ws.Range("L77").Value = 6714
$ws.Range("M77").Value = -3444.273
$ws.Range("N77").Value = -15450

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H103").Value = 25193.229
$ws.Range("J103").Value = 25193.229
$ws.Range("L103").Value = 25193.229
$ws.Range("N103").Value = -27537.229

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2198.158
$ws.Range("I122").Value = 1762.6471
$ws.Range("K122").Value = 5287.9413
$ws.Range("M122").Value = -2837.9413

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 253
$ws.Range("I4").Value = 221.42857
$ws.Range("J4").Value = 326.66666
$ws.Range("K4").Value = 221.42857
$ws.Range("L4").Value = 326.66666
$ws.Range("M4").Value = -106.42857
$ws.Range("N4").Value = -556.66666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 605.55
$ws.Range("I80").Value = 726.8889
$ws.Range("J80").Value = 506.27274
$ws.Range("K80").Value = 726.8889
$ws.Range("L80").Value = 506.27274
$ws.Range("M80").Value = 271.1111
$ws.Range("N80").Value = -2502.27274

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 605.55
$ws.Range("I83").Value = 726.8889
$ws.Range("J83").Value = 506.27274
$ws.Range("K83").Value = 3634.4445
$ws.Range("L83").Value = 2531.3637
$ws.Range("M83").Value = 1357.5555
$ws.Range("N83").Value = -12515.3637

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 40000
$ws.Range("J122").Value = 40000
$ws.Range("L122").Value = 40000
$ws.Range("N122").Value = -49800

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2209.8
$ws.Range("I99").Value = 1864
$ws.Range("J99").Value = 2440.3333
$ws.Range("K99").Value = 1864
$ws.Range("L99").Value = 2440.3333
$ws.Range("M99").Value = -366
$ws.Range("N99").Value = -5436.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2209.8
$ws.Range("I126").Value = 1864
$ws.Range("J126").Value = 2440.3333
$ws.Range("K126").Value = 5592
$ws.Range("L126").Value = 7320.999899999999
$ws.Range("M126").Value = -3122
$ws.Range("N126").Value = -12260.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 767.46155
$ws.Range("I5").Value = 337
$ws.Range("J5").Value = 1580.5555
$ws.Range("K5").Value = 1011
$ws.Range("L5").Value = 4741.666499999999
$ws.Range("M5").Value = -899
$ws.Range("N5").Value = -4965.666499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1436.8379
$ws.Range("J131").Value = 1204.4688
$ws.Range("L131").Value = 3613.4064
$ws.Range("N131").Value = -13693.4064

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 767.46155
$ws.Range("I135").Value = 337
$ws.Range("J135").Value = 1580.5555
$ws.Range("K135").Value = 3033
$ws.Range("L135").Value = 14224.9995
$ws.Range("M135").Value = -498
$ws.Range("N135").Value = -19294.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 44.31579
$ws.Range("I2").Value = 62.57143
$ws.Range("J2").Value = 33.666668
$ws.Range("K2").Value = 62.57143
$ws.Range("L2").Value = 33.666668
$ws.Range("M2").Value = 50.42857
$ws.Range("N2").Value = -259.666668

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 15000
$ws.Range("I122").Value = 20000
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 60000
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -57550
$ws.Range("N122").Value = -34900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4666.6665
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 6250
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 18750
$ws.Range("M126").Value = -2030
$ws.Range("N126").Value = -23690

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5265415
$ws.Range("J126").Value = 14288687
$ws.Range("L126").Value = 42866061
$ws.Range("N126").Value = -42871001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 363352.5
$ws.Range("I132").Value = 1112903.1
$ws.Range("J132").Value = 8302.210999999999
$ws.Range("K132").Value = 3338709.3
$ws.Range("L132").Value = 24906.633
$ws.Range("M132").Value = -3336179.3
$ws.Range("N132").Value = -29966.633
